# ---------------------------------------------------------------------------
# Edit script: 
#   1) Re-style the three tables (slides 14, 15, 16) from the custom
#      "Table_0" style {77D016E4-5A22-4E1C-AD43-6BFB0283C047} to the built-in
#      table style {045E1EA8-07B8-4261-9DCD-0DBDFBFB7E34}.
#   2) Swap the presentation's colour theme: the live "Integral" / "Red
#      Violet" colour scheme becomes the stock Office colour scheme
#      (black/white/blue-grey "Office" palette). Font scheme and format
#      scheme are already identical between the two theme parts, so only
#      the 12 theme colours need updating.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------------
$newTableStyle = "{045E1EA8-07B8-4261-9DCD-0DBDFBFB7E34}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyle)
    }
}

# --- 2) Theme colour scheme --------------------------------------------------
# Target ("Office") palette, as RGB hex (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink). ColorFormat.RGB takes a BGR-packed integer (standard
# Windows COLORREF), so convert each hex triplet accordingly.
$officeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$slide1 = $p.Slides.Item(1)
$colorScheme = $slide1.ThemeColorScheme

for ($i = 0; $i -lt $officeHex.Count; $i++) {
    $hex = $officeHex[$i]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $bgr = ($b * 65536) + ($g * 256) + $r
    $colorScheme.Item($i + 1).RGB = $bgr
}
